# Generate Report for Handoff
# Updates the "Latest Handoff Datetime" for the file
# 63373375-e236-49d5-ad70-28fdb75af732.md on both the zh-cn and de-de
# localization sheets (row 5, column H), and refreshes the corresponding
# "Latest HO Xliff Generate Date" on the Overview sheet (row 5, column G)
# to the newest of the two handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsZhCn.Range("H5").Value = "2016-10-18 10:36:05"
$wsDeDe.Range("H5").Value = "2016-10-18 10:36:26"
$wsOverview.Range("G5").Value = "2016-10-18 10:36:26"
